$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-23 (columns A-F), reflecting the reordering of
# controller rows as described by the commit ("Supported controller
# into new json data"). Rows 1, 24, 25, 26 are unchanged.
$data = @(
    @(902, 1, 0, 0, 0, 0),
    @(601, 9, 60, 67, 60, 42),
    @(201, 9, 30, 15, 45, 30),
    @(801, 3, 67, 65, 52, 45),
    @(1202, 2, 10, 10, 10, 10),
    @(1203, 3, 15, 15, 15, 15),
    @(901, 16, 15, 45, 60, 60),
    @(301, 6, 45, 30, 60, 45),
    @(701, 3, 90, 45, 97, 15),
    @(1201, 2, 10, 10, 10, 10),
    @(101, 9, 30, 15, 60, 15),
    @(1001, 18, 30, 75, 60, 72),
    @(501, 9, 52, 30, 75, 45),
    @(401, 9, 48, 67, 75, 45),
    @(1101, 0, 15, 30, 30, 0),
    @(502, 0, 4, 0, 0, 0),
    @(802, 0, 4, 5, 4, 0),
    @(3, 0, 3, 3, 3, 3),
    @(2, 0, 2, 2, 2, 2),
    @(1, 0, 2, 2, 2, 2),
    @(602, 0, 0, 4, 0, 9),
    @(402, 0, 0, 4, 0, 0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
